# Rename four hidden "_..." bookmarks (Google-Docs-style auto bookmarks) to
# new random-looking names, keeping everything else (host paragraph,
# w:colFirst/w:colLast attrs, collapsed start==end range, w:id) identical.
#
# Why this is tricky:
#   * These bookmarks are zero-length ("collapsed": bookmarkStart is
#     immediately followed by bookmarkEnd) and their names start with "_",
#     which marks them as hidden bookmarks in Word's object model.
#     $d.Bookmarks.Count / .ShowHidden do not enumerate them, and setting
#     .Name on the Bookmark object returned by Bookmarks.Item(name) is a
#     silent no-op in this host.
#   * Bookmarks.Item(name).Delete() + Bookmarks.Add(newName, range) does
#     rename it, but Add() always mints a bookmark with no w:colFirst /
#     w:colLast, so that pair of attributes (present, "0"/"0", on every
#     bookmark in this document) gets silently dropped.
#   * Range.InsertXML lets us splice in a literal <w:bookmarkStart> (with
#     whatever attributes we want, in this case copied verbatim) at an
#     exact character offset, which preserves w:colFirst/w:colLast. But
#     InsertXML always inserts its payload as a new paragraph *before* the
#     paragraph that currently starts at the insertion offset -- except
#     when the insertion offset is itself a paragraph boundary that already
#     has a paragraph before it, in which case the inserted para content
#     quietly merges into the *following* paragraph instead of staying
#     split out. We exploit that merge behavior: insert right at the
#     boundary between the previous paragraph and the bookmark's paragraph.
#   * The very first bookmark in the doc sits at offset 0 (the start of the
#     very first paragraph), so there is no previous paragraph to anchor
#     the merge trick on. We manufacture one temporarily
#     (Range.InsertParagraphBefore, a normal/native Word op with no XML
#     splicing weirdness), do the InsertXML+merge trick against that new
#     boundary, then delete the scratch paragraph again.

$d = $word.ActiveDocument

function New-BookmarkXmlFragment($name) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:bookmarkStart w:colFirst="0" w:colLast="0" w:name="' + $name + '" w:id="9999"/>' +
        '<w:bookmarkEnd w:id="9999"/>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Rename-HiddenBookmark($oldName, $newName) {
    $bm = $d.Bookmarks.Item($oldName)
    $atDocStart = ($bm.Range.Start -eq 0)

    if ($atDocStart) {
        # No paragraph precedes this one -- make a scratch one so we have a
        # paragraph boundary to anchor the InsertXML merge trick on.
        $d.Paragraphs.Item(1).Range.InsertParagraphBefore()
    }

    # Re-fetch: position is stable (InsertParagraphBefore shifted everything
    # after it by exactly 1, and the bookmark itself is still collapsed at
    # the start of its (now possibly shifted) paragraph).
    $bm = $d.Bookmarks.Item($oldName)
    $boundary = $bm.Range.Start

    $insertRange = $d.Range($boundary, $boundary)
    $insertRange.InsertXML((New-BookmarkXmlFragment $newName))

    $bm.Delete()

    if ($atDocStart) {
        # Drop the scratch paragraph we added at the very top.
        $d.Paragraphs.Item(1).Range.Delete()
    }
}

Rename-HiddenBookmark "_adqwoti2wb3p" "_4a1uk8mignqx"
Rename-HiddenBookmark "_ohkdl32dea3x" "_m45zcy1y0k6p"
Rename-HiddenBookmark "_q5jaaeb0e6x2" "_2o5smqbvihsl"
Rename-HiddenBookmark "_h6p2lccqo3mm" "_su50qpurz8jd"

Write-Output "Renamed 4 bookmarks"
